# Insert a new weekly price record as row 401, pushing the existing
# rows 401:421 down to 402:422 (dimension grows from A1:R421 to A1:R422).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 401..421 down by one row.
$ws.Rows.Item(401).Insert()

# Populate the newly inserted row 401 with the new record.
$ws.Cells.Item(401, 1).Value = 11
$ws.Cells.Item(401, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(401, 3).Value = "Bíobío"
$ws.Cells.Item(401, 4).Value = 45008
$ws.Cells.Item(401, 5).Value = 8
$ws.Cells.Item(401, 6).Value = 100114001
$ws.Cells.Item(401, 7).Value = "Papa"
$ws.Cells.Item(401, 8).Value = "Asterix"
$ws.Cells.Item(401, 9).Value = "1a (cosecha)"
$ws.Cells.Item(401, 10).Value = 5000
$ws.Cells.Item(401, 11).Value = 11500
$ws.Cells.Item(401, 12).Value = 12000
$ws.Cells.Item(401, 13).Value = 11750
$ws.Cells.Item(401, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(401, 15).Value = "Provincia de Arauco"
$ws.Cells.Item(401, 16).Value = 470
$ws.Cells.Item(401, 17).Value = 25
$ws.Cells.Item(401, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Range("D401").NumberFormat = $ws.Range("D402").NumberFormat
